# Apply the corrections described in the commit diff:
# 1. Fix typo in category label: "השקה - משכנתא" -> "השקעה - משכנתא"
# 2. Update the sheet view's scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("initial_categories")

# 1. Fix the typo in cell A43 (shared string used there)
$ws.Range("A43").Value = "השקעה - משכנתא"

# 2. Update the visible top-left cell and the current selection
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B37:B38").Select()
